$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (shared strings) ---
$ws.Range("B1").Value = "TISG"
$ws.Range("D1").Value = "buy"
$ws.Range("E1").Value = "MYDIR"
$ws.Range("F1").Value = "need_to_buy"

# --- Data updates (columns B, C, E, F for rows 2-15; A and D unchanged) ---
$ws.Range("B2").Value = 7363.23288217521
$ws.Range("C2").Value = 7084.09096009052
$ws.Range("E2").Value = 3212.56170555259
$ws.Range("F2").Value = -1.80613893153698

$ws.Range("B3").Value = 7372.65843003816
$ws.Range("C3").Value = 7036.07286381342
$ws.Range("E3").Value = 3223.42373793971
$ws.Range("F3").Value = 276.645691739714

$ws.Range("B4").Value = 6956.53411314305
$ws.Range("C4").Value = 6645.51861749365
$ws.Range("E4").Value = 3127.38889593348
$ws.Range("F4").Value = 256.371146392797

$ws.Range("B5").Value = 6657.92668775434
$ws.Range("C5").Value = 6373.40814371494
$ws.Range("E5").Value = 2905.87211299751
$ws.Range("F5").Value = 235.803344029685

$ws.Range("B6").Value = 6590.10304858047
$ws.Range("C6").Value = 5850.3081508992
$ws.Range("E6").Value = 2878.33462844333
$ws.Range("F6").Value = 212.860115805939

$ws.Range("B7").Value = 2211.54596239233
$ws.Range("C7").Value = 4269.58048343521
$ws.Range("E7").Value = 2714.01002757637
$ws.Range("F7").Value = 140.149604625483

$ws.Range("B8").Value = 2228.93158091655
$ws.Range("C8").Value = 4253.83126795722
$ws.Range("E8").Value = 2821.51964949253
$ws.Range("F8").Value = 143.972954893739

$ws.Range("B9").Value = 7225.19876436481
$ws.Range("C9").Value = 7489.06635504889
$ws.Range("E9").Value = 3511.35451968222
$ws.Range("F9").Value = 307.51753644713

$ws.Range("B10").Value = 7225.19876436481
$ws.Range("C10").Value = 8014.79199974509
$ws.Range("E10").Value = 3511.35451968222
$ws.Range("F10").Value = 329.422771642805

$ws.Range("B11").Value = 7225.19876436481
$ws.Range("C11").Value = 7710.72090179961
$ws.Range("E11").Value = 3511.35451968222
$ws.Range("F11").Value = 316.753142561743

$ws.Range("B12").Value = 7225.19876436481
$ws.Range("C12").Value = 7623.7412200405
$ws.Range("E12").Value = 3511.35451968222
$ws.Range("F12").Value = 313.128989155113

$ws.Range("B13").Value = 6107.6829751599
$ws.Range("C13").Value = 6183.40699096311
$ws.Range("E13").Value = 3446.71757656699
$ws.Range("F13").Value = 250.421856980421

$ws.Range("B14").Value = 2557.17780844876
$ws.Range("C14").Value = 4328.79164727458
$ws.Range("E14").Value = 3796.69318447108
$ws.Range("F14").Value = 61.3952013227359

$ws.Range("B15").Value = 2557.17780844876
$ws.Range("C15").Value = 4353.1497113995
$ws.Range("E15").Value = 3796.69318447108
$ws.Range("F15").Value = 62.4101206612744
